# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1930.5264
$ws.Range("J40").Value = 1854.5555
$ws.Range("L40").Value = 1854.5555
$ws.Range("N40").Value = -2204.5555

# Row 64
$ws.Range("H64").Value = 3466.8987
$ws.Range("I64").Value = 3538.0334
$ws.Range("J64").Value = 3242.2632
$ws.Range("K64").Value = 3538.0334
$ws.Range("L64").Value = 3242.2632
$ws.Range("M64").Value = -3290.0334
$ws.Range("N64").Value = -3738.2632

# Row 67
$ws.Range("H67").Value = 3466.8987
$ws.Range("I67").Value = 3538.0334
$ws.Range("J67").Value = 3242.2632
$ws.Range("K67").Value = 3538.0334
$ws.Range("L67").Value = 3242.2632
$ws.Range("M67").Value = -2680.0334
$ws.Range("N67").Value = -4958.263199999999

# Row 98
$ws.Range("H98").Value = 1546.6786
$ws.Range("I98").Value = 1620.65
$ws.Range("K98").Value = 1620.65
$ws.Range("M98").Value = -122.6500000000001

# Row 100
$ws.Range("H100").Value = 11907119
$ws.Range("I100").Value = 18520520
$ws.Range("K100").Value = 18520520
$ws.Range("M100").Value = -18519979

# Row 116
$ws.Range("H116").Value = 6909.8184
$ws.Range("I116").Value = 8622.143
$ws.Range("J116").Value = 3913.25
$ws.Range("K116").Value = 8622.143
$ws.Range("L116").Value = 3913.25
$ws.Range("M116").Value = -5180.143
$ws.Range("N116").Value = -10797.25

# Row 122
$ws.Range("H122").Value = 1546.6786
$ws.Range("I122").Value = 1620.65
$ws.Range("K122").Value = 4861.950000000001
$ws.Range("M122").Value = -2411.950000000001

# Row 132
$ws.Range("H132").Value = 2703.2144
$ws.Range("I132").Value = 2887.0833
$ws.Range("J132").Value = 1600
$ws.Range("K132").Value = 8661.249899999999
$ws.Range("L132").Value = 4800
$ws.Range("M132").Value = -6131.249899999999
$ws.Range("N132").Value = -9860

# Row 137
$ws.Range("H137").Value = 1719.4445
$ws.Range("I137").Value = 1646.7858
$ws.Range("J137").Value = 1839.1177
$ws.Range("K137").Value = 4940.357400000001
$ws.Range("L137").Value = 5517.3531
$ws.Range("M137").Value = -2390.357400000001
$ws.Range("N137").Value = -10617.3531

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3276.4285
$ws.Range("I2").Value = 2540.6667
$ws.Range("J2").Value = 3828.25
$ws.Range("K2").Value = 2540.6667
$ws.Range("L2").Value = 3828.25
$ws.Range("M2").Value = -2427.6667
$ws.Range("N2").Value = -4054.25

# Row 61
$ws.Range("H61").Value = 6751
$ws.Range("I61").Value = 7580.7646
$ws.Range("J61").Value = 4400
$ws.Range("K61").Value = 7580.7646
$ws.Range("L61").Value = 4400
$ws.Range("M61").Value = -7368.7646
$ws.Range("N61").Value = -4824

# Row 116
$ws.Range("H116").Value = 3276.4285
$ws.Range("I116").Value = 2540.6667
$ws.Range("J116").Value = 3828.25
$ws.Range("K116").Value = 2540.6667
$ws.Range("L116").Value = 3828.25
$ws.Range("M116").Value = -246.6667000000002
$ws.Range("N116").Value = -8416.25

# Row 136
$ws.Range("H136").Value = 6751
$ws.Range("I136").Value = 7580.7646
$ws.Range("J136").Value = 4400
$ws.Range("K136").Value = 22742.2938
$ws.Range("L136").Value = 13200
$ws.Range("M136").Value = -20192.2938
$ws.Range("N136").Value = -18300

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3276.4285
$ws.Range("I3").Value = 2540.6667
$ws.Range("J3").Value = 3828.25
$ws.Range("K3").Value = 2540.6667
$ws.Range("L3").Value = 3828.25
$ws.Range("M3").Value = -2426.6667
$ws.Range("N3").Value = -4056.25

# Row 134
$ws.Range("H134").Value = 3741.6304
$ws.Range("I134").Value = 3832.0732
$ws.Range("K134").Value = 11496.2196
$ws.Range("M134").Value = -8961.2196

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 19745
$ws.Range("I31").Value = 1806.3529
$ws.Range("J31").Value = 40075.465
$ws.Range("K31").Value = 1806.3529
$ws.Range("L31").Value = 40075.465
$ws.Range("M31").Value = -1511.3529
$ws.Range("N31").Value = -40665.465

# Row 34
$ws.Range("H34").Value = 19745
$ws.Range("I34").Value = 1806.3529
$ws.Range("J34").Value = 40075.465
$ws.Range("K34").Value = 1806.3529
$ws.Range("L34").Value = 40075.465
$ws.Range("M34").Value = -1604.3529
$ws.Range("N34").Value = -40479.465

# Row 55
$ws.Range("H55").Value = 13277.5
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 13277.5
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 13277.5
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -13907.5

# Row 58
$ws.Range("H58").Value = 2158.8857
$ws.Range("I58").Value = 1427.3572
$ws.Range("J58").Value = 2646.5715
$ws.Range("K58").Value = 1427.3572
$ws.Range("L58").Value = 2646.5715
$ws.Range("M58").Value = -1224.3572
$ws.Range("N58").Value = -3052.5715

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 132
$ws.Range("H132").Value = 2612.1765
$ws.Range("I132").Value = 2136.4546
$ws.Range("J132").Value = 3484.3333
$ws.Range("K132").Value = 6409.3638
$ws.Range("L132").Value = 10452.9999
$ws.Range("M132").Value = -3879.3638
$ws.Range("N132").Value = -15512.9999

# Row 134
$ws.Range("H134").Value = 3999.0435
$ws.Range("I134").Value = 4285.4443
$ws.Range("J134").Value = 2968
$ws.Range("K134").Value = 12856.3329
$ws.Range("L134").Value = 8904
$ws.Range("M134").Value = -10321.3329
$ws.Range("N134").Value = -13974

# Row 136
$ws.Range("H136").Value = 2158.8857
$ws.Range("I136").Value = 1427.3572
$ws.Range("J136").Value = 2646.5715
$ws.Range("K136").Value = 4282.071599999999
$ws.Range("L136").Value = 7939.7145
$ws.Range("M136").Value = -1732.071599999999
$ws.Range("N136").Value = -13039.7145

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 239134.5
$ws.Range("I5").Value = 546.6667
$ws.Range("J5").Value = 334569.62
$ws.Range("K5").Value = 1640.0001
$ws.Range("L5").Value = 1003708.86
$ws.Range("M5").Value = -1528.0001
$ws.Range("N5").Value = -1003932.86

# Row 132
$ws.Range("H132").Value = 2258.5264
$ws.Range("I132").Value = 1491.8182
$ws.Range("J132").Value = 3312.75
$ws.Range("K132").Value = 13426.3638
$ws.Range("L132").Value = 29814.75
$ws.Range("M132").Value = -10896.3638
$ws.Range("N132").Value = -34874.75

# Row 135
$ws.Range("H135").Value = 239134.5
$ws.Range("I135").Value = 546.6667
$ws.Range("J135").Value = 334569.62
$ws.Range("K135").Value = 4920.0003
$ws.Range("L135").Value = 3011126.58
$ws.Range("M135").Value = -2385.0003
$ws.Range("N135").Value = -3016196.58

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 4380.1704
$ws.Range("I126").Value = 6750.1055
$ws.Range("J126").Value = 2772
$ws.Range("K126").Value = 20250.3165
$ws.Range("L126").Value = 8316
$ws.Range("M126").Value = -17780.3165
$ws.Range("N126").Value = -13256

# Row 132
$ws.Range("H132").Value = 2521.48
$ws.Range("I132").Value = 1945.9231
$ws.Range("K132").Value = 5837.7693
$ws.Range("M132").Value = -3307.7693

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1740.9166
$ws.Range("I16").Value = 1740.9166
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1740.9166
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1570.9166
$ws.Range("N16").ClearContents()

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 93
$ws.Range("H93").Value = 1450
$ws.Range("I93").Value = 1400
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 1400
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = -152
$ws.Range("N93").Value = -3996

# Row 132
$ws.Range("H132").Value = 15879996
$ws.Range("I132").Value = 25649840
$ws.Range("K132").Value = 76949520
$ws.Range("M132").Value = -76946990

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 63129.938
$ws.Range("I100").Value = 125423.625
$ws.Range("K100").Value = 250847.25
$ws.Range("M100").Value = -250306.25
